$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$labelCol = $ws.Range("A1:A22")

# Version: 1.8.1 -> 1.8.2
$row = $labelCol.Find("Version").Row
$ws.Cells.Item($row, 2).Value = "1.8.2"

# Status: draft -> active
$row = $labelCol.Find("Status").Row
$ws.Cells.Item($row, 2).Value = "active"

# Experimental: true -> (cleared, CodeSystem is no longer experimental)
$row = $labelCol.Find("Experimental").Row
$ws.Cells.Item($row, 2).ClearContents()

# Date: 2024-01-18 -> 2025-11-18 (kept as text, not auto-converted to a date serial)
$row = $labelCol.Find("Date").Row
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2025-11-18"
